$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 708.56525
$ws.Range("I28").Value = 297.91666
$ws.Range("J28").Value = 1156.5454
$ws.Range("K28").Value = 297.91666
$ws.Range("L28").Value = 1156.5454
$ws.Range("M28").Value = 187.08334
$ws.Range("N28").Value = -2126.5454
$ws.Range("H32").Value = 1000
$ws.Range("J32").Value = 1000
$ws.Range("L32").Value = 1000
$ws.Range("N32").Value = -1652
$ws.Range("H62").Value = 6829.769
$ws.Range("I62").Value = 4842.143
$ws.Range("J62").Value = 9148.666999999999
$ws.Range("K62").Value = 4842.143
$ws.Range("L62").Value = 9148.666999999999
$ws.Range("M62").Value = -4218.143
$ws.Range("N62").Value = -10396.667
$ws.Range("H65").Value = 6829.769
$ws.Range("I65").Value = 4842.143
$ws.Range("J65").Value = 9148.666999999999
$ws.Range("K65").Value = 24210.715
$ws.Range("L65").Value = 45743.335
$ws.Range("M65").Value = -21090.715
$ws.Range("N65").Value = -51983.335
$ws.Range("H86").Value = 4479.6665
$ws.Range("I86").Value = 761.3570999999999
$ws.Range("J86").Value = 8484
$ws.Range("K86").Value = 761.3570999999999
$ws.Range("L86").Value = 8484
$ws.Range("M86").Value = 361.6429000000001
$ws.Range("N86").Value = -10730
$ws.Range("H88").Value = 587.4
$ws.Range("I88").Value = 193.33333
$ws.Range("J88").Value = 756.2857
$ws.Range("K88").Value = 193.33333
$ws.Range("L88").Value = 756.2857
$ws.Range("M88").Value = 212.66667
$ws.Range("N88").Value = -1568.2857
$ws.Range("H89").Value = 4479.6665
$ws.Range("I89").Value = 761.3570999999999
$ws.Range("J89").Value = 8484
$ws.Range("K89").Value = 3806.7855
$ws.Range("L89").Value = 42420
$ws.Range("M89").Value = 1809.2145
$ws.Range("N89").Value = -53652
$ws.Range("H91").Value = 587.4
$ws.Range("I91").Value = 193.33333
$ws.Range("J91").Value = 756.2857
$ws.Range("K91").Value = 193.33333
$ws.Range("L91").Value = 756.2857
$ws.Range("M91").Value = 1210.66667
$ws.Range("N91").Value = -3564.2857
$ws.Range("H98").Value = 761
$ws.Range("I98").Value = 761
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 761
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 737
$ws.Range("N98").ClearContents()
$ws.Range("H107").Value = 872
$ws.Range("I107").Value = 651.5238000000001
$ws.Range("K107").Value = 651.5238000000001
$ws.Range("M107").Value = 1268.4762
$ws.Range("H122").Value = 761
$ws.Range("I122").Value = 761
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2283
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 167
$ws.Range("N122").ClearContents()

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1466.579
$ws.Range("I2").Value = 1292.7646
$ws.Range("K2").Value = 1292.7646
$ws.Range("M2").Value = -1179.7646
$ws.Range("H32").Value = 7397.237
$ws.Range("I32").Value = 7117
$ws.Range("K32").Value = 7117
$ws.Range("M32").Value = -6830
$ws.Range("H45").Value = 2902.5625
$ws.Range("I45").Value = 3062.375
$ws.Range("J45").Value = 2742.75
$ws.Range("K45").Value = 3062.375
$ws.Range("L45").Value = 2742.75
$ws.Range("M45").Value = -2685.375
$ws.Range("N45").Value = -3496.75
$ws.Range("H61").Value = 7830.2856
$ws.Range("I61").Value = 7830.2856
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 7830.2856
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -7618.2856
$ws.Range("N61").ClearContents()
$ws.Range("H97").Value = 200000660
$ws.Range("I97").Value = 655
$ws.Range("K97").Value = 655
$ws.Range("M97").Value = -159
$ws.Range("H116").Value = 1466.579
$ws.Range("I116").Value = 1292.7646
$ws.Range("K116").Value = 1292.7646
$ws.Range("M116").Value = 1001.2354
$ws.Range("H136").Value = 7830.2856
$ws.Range("I136").Value = 7830.2856
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 23490.8568
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -20940.8568
$ws.Range("N136").ClearContents()

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1466.579
$ws.Range("I3").Value = 1292.7646
$ws.Range("K3").Value = 1292.7646
$ws.Range("M3").Value = -1178.7646
$ws.Range("H86").Value = 1502.54
$ws.Range("I86").Value = 1374.2333
$ws.Range("J86").Value = 1695
$ws.Range("K86").Value = 1374.2333
$ws.Range("L86").Value = 1695
$ws.Range("M86").Value = -251.2333000000001
$ws.Range("N86").Value = -3941
$ws.Range("H89").Value = 1502.54
$ws.Range("I89").Value = 1374.2333
$ws.Range("J89").Value = 1695
$ws.Range("K89").Value = 6871.1665
$ws.Range("L89").Value = 8475
$ws.Range("M89").Value = -1255.1665
$ws.Range("N89").Value = -19707

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1170
$ws.Range("I16").Value = 250
$ws.Range("K16").Value = 250
$ws.Range("M16").Value = 37
$ws.Range("H58").Value = 20238.223
$ws.Range("I58").Value = 1570.7142
$ws.Range("J58").Value = 40341.69
$ws.Range("K58").Value = 1570.7142
$ws.Range("L58").Value = 40341.69
$ws.Range("M58").Value = -1367.7142
$ws.Range("N58").Value = -40747.69
$ws.Range("H113").Value = 1170
$ws.Range("I113").Value = 250
$ws.Range("K113").Value = 250
$ws.Range("M113").Value = 1920
$ws.Range("H136").Value = 20238.223
$ws.Range("I136").Value = 1570.7142
$ws.Range("J136").Value = 40341.69
$ws.Range("K136").Value = 4712.142599999999
$ws.Range("L136").Value = 121025.07
$ws.Range("M136").Value = -2162.142599999999
$ws.Range("N136").Value = -126125.07

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1566.2106
$ws.Range("I5").Value = 1418.5
$ws.Range("J5").Value = 1673.6364
$ws.Range("K5").Value = 4255.5
$ws.Range("L5").Value = 5020.9092
$ws.Range("M5").Value = -4143.5
$ws.Range("N5").Value = -5244.9092
$ws.Range("H92").Value = 540.9286
$ws.Range("I92").Value = 266.66666
$ws.Range("J92").Value = 1034.6
$ws.Range("K92").Value = 799.9999799999999
$ws.Range("L92").Value = 3103.8
$ws.Range("M92").Value = 448.0000200000001
$ws.Range("N92").Value = -5599.799999999999
$ws.Range("H107").Value = 3897.926
$ws.Range("I107").Value = 6840
$ws.Range("J107").Value = 220.33333
$ws.Range("K107").Value = 20520
$ws.Range("L107").Value = 660.99999
$ws.Range("M107").Value = -18600
$ws.Range("N107").Value = -4500.99999
$ws.Range("H131").Value = 713.65
$ws.Range("J131").Value = 725.15955
$ws.Range("L131").Value = 2175.47865
$ws.Range("N131").Value = -12255.47865
$ws.Range("H135").Value = 1566.2106
$ws.Range("I135").Value = 1418.5
$ws.Range("J135").Value = 1673.6364
$ws.Range("K135").Value = 12766.5
$ws.Range("L135").Value = 15062.7276
$ws.Range("M135").Value = -10231.5
$ws.Range("N135").Value = -20132.7276
$ws.Range("H141").Value = 1885.7142
$ws.Range("I141").Value = 1885.7142
$ws.Range("K141").Value = 5657.142599999999
$ws.Range("M141").Value = -477.1425999999992

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 318.33334
$ws.Range("I107").Value = 318.875
$ws.Range("J107").Value = 317.25
$ws.Range("K107").Value = 318.875
$ws.Range("L107").Value = 317.25
$ws.Range("M107").Value = 1601.125
$ws.Range("N107").Value = -4157.25
$ws.Range("H126").Value = 5349.7417
$ws.Range("I126").Value = 4361.778
$ws.Range("J126").Value = 6717.6924
$ws.Range("K126").Value = 13085.334
$ws.Range("L126").Value = 20153.0772
$ws.Range("M126").Value = -10615.334
$ws.Range("N126").Value = -25093.0772
$ws.Range("H132").Value = 35883.6
$ws.Range("I132").Value = 2054
$ws.Range("K132").Value = 6162
$ws.Range("M132").Value = -3632

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2976.5
$ws.Range("I40").Value = 2353.7856
$ws.Range("K40").Value = 2353.7856
$ws.Range("M40").Value = -2217.7856
$ws.Range("H68").Value = 2611.7334
$ws.Range("I68").Value = 2470.9092
$ws.Range("J68").Value = 2999
$ws.Range("K68").Value = 2470.9092
$ws.Range("L68").Value = 2999
$ws.Range("M68").Value = -1721.9092
$ws.Range("N68").Value = -4497
$ws.Range("H71").Value = 2611.7334
$ws.Range("I71").Value = 2470.9092
$ws.Range("J71").Value = 2999
$ws.Range("K71").Value = 12354.546
$ws.Range("L71").Value = 14995
$ws.Range("M71").Value = -8610.546
$ws.Range("N71").Value = -22483
$ws.Range("H136").Value = 2382
$ws.Range("I136").Value = 2511
$ws.Range("K136").Value = 7533
$ws.Range("M136").Value = -4983

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 6494631.5
$ws.Range("I107").Value = 849.5
$ws.Range("K107").Value = 2548.5
$ws.Range("M107").Value = -628.5
